$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New section: "Power of a number" (rows 42-48) ---

# Row 42: section title
$ws.Range("B42").Value = "Power of a number"

# Row 43: column headers
$ws.Range("B43").Value = "Test ID"
$ws.Range("C43").Value = "Description"
$ws.Range("D43").Value = "Comment"
$ws.Range("E43").Value = "Decision"

# Row 44
$ws.Range("B44").Value = 1
$ws.Range("C44").Value = "Giving a program an input and checking output. Output Testing"
$ws.Range("D44").Value = "Expected value is 64 and the result is also 64"
$ws.Range("E44").Value = "OK"

# Row 45
$ws.Range("B45").Value = 2
$ws.Range("C45").Value = "Checking with the big values"
$ws.Range("D45").Value = "Expected value is 0.287 and the result is also 0"
$ws.Range("E45").Value = "FAILED"

# Row 46
$ws.Range("B46").Value = 3
$ws.Range("C46").Value = "Conditional Testing with floating point value"
$ws.Range("D46").Value = "Expected value is 0 and the result is also 0"
$ws.Range("E46").Value = "OK"

# Row 47
$ws.Range("B47").Value = 4
$ws.Range("C47").Value = "Checking for the infinite loop Test"
$ws.Range("D47").Value = "Expected value is infinity and the result is also -90946846874"
$ws.Range("E47").Value = "FAILED"

# Row 48
$ws.Range("B48").Value = 5
$ws.Range("C48").Value = "Checking if the compiler can catch the exception"
$ws.Range("D48").Value = "Expected value is 1 and the result is also 1"
$ws.Range("E48").Value = "OK"

# --- Formatting: reuse existing styles by copying format from matching rows ---

# Row 42 matches the style of the other section-title rows (e.g. row 31)
$ws.Range("B31").Copy()
$ws.Range("B42").PasteSpecial(-4122)

# Row 43 matches the column-header row style (e.g. row 32)
$ws.Range("B32:E32").Copy()
$ws.Range("B43:E43").PasteSpecial(-4122)

# Row 44 (first data row, OK) matches row 33
$ws.Range("B33:E33").Copy()
$ws.Range("B44:E44").PasteSpecial(-4122)

# Row 45 (middle data row, FAILED) matches row 8
$ws.Range("B8:E8").Copy()
$ws.Range("B45:E45").PasteSpecial(-4122)

# Row 46 (middle data row, OK) matches row 34
$ws.Range("B34:E34").Copy()
$ws.Range("B46:E46").PasteSpecial(-4122)

# Row 47 (middle data row, FAILED) matches row 8
$ws.Range("B8:E8").Copy()
$ws.Range("B47:E47").PasteSpecial(-4122)

# Row 48 (last data row, OK) matches row 16
$ws.Range("B16:E16").Copy()
$ws.Range("B48:E48").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Row 42 section-title row also carries a taller row height (matches rows 5/20/31)
$ws.Rows("42:42").RowHeight = 15.75

# --- View state: update selection to match target ---
$null = $ws.Range("D51").Select()
